$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$r = $d.Range($p1.Range.Start, $p1.Range.Start)
$r.InsertAfter("Z")
Write-Output "done"
